$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset the header row's formatting (drops the pink fill style) and
# wipe the old sample-data rows (2-4) completely.
$ws.Rows.Item(1).ClearFormats()
$ws.Range("A2:C4").Clear()

# Write the new header row values
$ws.Range("A1").Value = "Employee ID"
$ws.Range("B1").Value = "First and Middle Name"
$ws.Range("C1").Value = "Last Name"

# Set column widths to match the target layout
$ws.Columns.Item(1).ColumnWidth = 23.44140625
$ws.Columns.Item(2).ColumnWidth = 29.88671875
$ws.Columns.Item(3).ColumnWidth = 27

# Move the active selection
$ws.Range("E17").Select()
